# Edit script: add "2022-Q1" sheet with fund holding data, and
# insert a corresponding summary row at the top of the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q1" sheet right after "2021-Q4" (position 5),
#    so the tab order becomes: ... 2021-Q4, 2022-Q1, 总计
# ---------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item(5)
$q4Sheet = $wb.Worksheets.Item(5)
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "2022-Q1"

# Copy the bold/bordered header style (used on row 1, columns B:H) and the
# index-column style (used on column A of data rows) from the existing
# "2021-Q4" sheet, which has the identical 8-column layout.
$q4Sheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$q4Sheet.Range("A2").Copy()
$newSheet.Range("A2:A29").PasteSpecial(-4122)

# Header row
$headers = @('基金代码', '基金名称', '基金规模', '股票总仓位', '仓位占比', '持有市值(亿元)', '仓位排名')
for ($c = 0; $c -lt $headers.Length; $c++) {
    $newSheet.Cells.Item(1, $c + 2).Value = $headers[$c]
}

# Force columns B:G (fund code/name/scale/position/ratio/market-value) to be
# stored as text, matching the source data (keeps leading zeros in fund
# codes and keeps the decimal strings as literal text rather than floats).
$newSheet.Range("B2:G29").NumberFormat = "@"

$fundRows = @(
    @(0, '002340', '富国价值优势混合', '69.56', '93.66', '2.98', '2.0729', 6),
    @(1, '011578', '汇丰晋信核心成长混合型证券投资基金A', '31.02', '91.97', '4.53', '1.4052', 7),
    @(2, '001643', '汇丰晋信智造先锋股票A', '29.09', '92.99', '4.59', '1.3352', 5),
    @(3, '010400', '富国融泰三个月定期开放混合', '34.05', '94.54', '2.81', '0.9568', 5),
    @(4, '004674', '富国新机遇灵活配置混合A', '23.28', '93.61', '2.57', '0.5983', 7),
    @(5, '001644', '汇丰晋信智造先锋股票C', '10.91', '92.99', '4.59', '0.5008', 5),
    @(6, '009334', '富国融享18个月定期开放混合', '15.00', '94.12', '2.86', '0.4290', 6),
    @(7, '001113', '南方大数据100指数A', '20.79', '94.23', '1.94', '0.4033', 5),
    @(8, '002124', '广发新兴产业精选灵活配置混合', '11.11', '91.11', '2.44', '0.2711', 10),
    @(9, '011579', '汇丰晋信核心成长混合型证券投资基金C', '4.27', '91.97', '4.53', '0.1934', 7),
    @(10, '004604', '富国新活力灵活配置混合A', '5.72', '94.03', '3.24', '0.1853', 9),
    @(11, '161039', '富国中证1000指数增强LOF', '21.72', '89.03', '0.83', '0.1803', 3),
    @(12, '011410', '中信建投量化进取6个月持有期混合A', '9.13', '93.80', '1.20', '0.1096', 1),
    @(13, '004675', '富国新机遇灵活配置混合C', '3.84', '93.61', '2.57', '0.0987', 7),
    @(14, '012878', '中信建投量化精选6个月持有期混合型证券投资基金A', '6.32', '88.35', '0.86', '0.0544', 9),
    @(15, '004605', '富国新活力灵活配置混合C', '1.41', '94.03', '3.24', '0.0457', 9),
    @(16, '012879', '中信建投量化精选6个月持有期混合型证券投资基金C', '4.26', '88.35', '0.86', '0.0366', 9),
    @(17, '011411', '中信建投量化进取6个月持有期混合C', '2.15', '93.80', '1.20', '0.0258', 1),
    @(18, '006181', '格林伯锐灵活配置混合A', '0.29', '89.68', '6.56', '0.0190', 1),
    @(19, '004194', '招商中证1000指数增强A', '1.76', '94.40', '1.05', '0.0185', 8),
    @(20, '013242', '北信瑞丰优势行业股票', '0.82', '92.63', '1.17', '0.0096', 6),
    @(21, '006182', '格林伯锐灵活配置混合C', '0.12', '89.68', '6.56', '0.0079', 1),
    @(22, '004195', '招商中证1000指数增强C', '0.68', '94.40', '1.05', '0.0071', 8),
    @(23, '005536', '渤海汇金量化成长混合', '0.61', '88.57', '0.82', '0.0050', 4),
    @(24, '008300', '人保量化锐进混合A', '0.20', '92.95', '2.27', '0.0045', 9),
    @(25, '004344', '南方大数据100指数C', '0.17', '94.23', '1.94', '0.0033', 5),
    @(26, '002952', '建信多因子量化股票', '0.10', '91.47', '3.10', '0.0031', 3),
    @(27, '008301', '人保量化锐进混合C', '0.06', '92.95', '2.27', '0.0014', 9)
)

foreach ($row in $fundRows) {
    $r = [int]$row[0] + 2
    $newSheet.Cells.Item($r, 1).Value = [int]$row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $newSheet.Cells.Item($r, 5).Value = $row[4]
    $newSheet.Cells.Item($r, 6).Value = $row[5]
    $newSheet.Cells.Item($r, 7).Value = $row[6]
    $newSheet.Cells.Item($r, 8).Value = [int]$row[7]
}

# ---------------------------------------------------------------------
# 2) Insert a new row at the top of the "总计" sheet's data (row 2) for
#    the 2022-Q1 summary, pushing the existing quarters down by one row.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("A2:D2").ClearFormats()

# Re-apply the bold index-column style to A2 (copied from A3, which still
# carries the original style after the insert shifted it down).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 28
$totalSheet.Range("D2").Value = 8.98

# Re-number the index column (A) sequentially 0..5 for the 6 data rows
# now that the new row has shifted everything else down by one.
for ($r = 3; $r -le 7; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}

Write-Host "2022-Q1 sheet and 总计 summary row added."
